# Refresh the cryptos price list (GitHub Actions scheduled update).
# Price cells in column D (e.g. "7.62", "26.10") look like numbers, so a
# plain .Value assignment would make Excel auto-convert them to the Number
# type (dropping trailing zeros / the text formatting). Force those cells
# to Text first, write the value, then reset the style back to "Normal" so
# no stray number-format style is left attached to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.930.00"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.18%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.232.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.75%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "539.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.89%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.63"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.98%  "

$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.233.90"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.66%  "

$ws.Range("E9").Value = "  -4.25%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.62"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.99%  "

$ws.Range("E11").Value = "  -5.62%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.396"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.10%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.788.37"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.70%  "

$ws.Range("E14").Value = "  -1.23%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.10"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.75%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.231.56"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.98%  "

$ws.Range("E17").Value = "  -5.79%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "58.978.28"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.27%  "

$ws.Range("E19").Value = "  -7.03%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.31"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.18%  "

$ws.Range("E21").Value = "  -6.22%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "361.86"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.20%  "

$ws.Range("E23").Value = "  -0.10%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.60"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.38%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.521"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -6.89%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.367.52"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.88%  "

$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.170"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.04%  "

$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0976"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -10.20%  "

$ws.Range("E29").Value = "  +0.17%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.10"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.62%  "

$ws.Range("E31").Value = "  -0.03%  "

$ws.Range("E32").Value = "  -6.51%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.11"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.22%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "21.97"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.88%  "

$ws.Range("E35").Value = "  -4.60%  "

$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.95"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.57%  "

$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "163.12"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.24%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.44"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.84%  "

$ws.Range("E39").Value = "  -6.71%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "26.60"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -9.06%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0714"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.92%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.266.94"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.80%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.17"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.67%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.717"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.63%  "

$ws.Range("E45").Value = "  -3.41%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.04"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.83%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.51"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.26%  "

$ws.Range("E48").Value = "  -0.05%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.297.04"
$ws.Range("D49").Style = "Normal"

$ws.Range("E50").Value = "  -5.66%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.94"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.40%  "
